# Apply "Optuna Attempt (go back with original)" edits.
# Sheet "Forecast Comparison": update MyForecast (D) and Seasonality Index (L) values.
# Sheet "Summary": update forecast total values (B9:B12).

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: column L (Seasonality Index) ---
$seasonality = @{
    2  = 1.2
    3  = 1.16
    4  = 1.16
    5  = 1.11
    6  = 0.9399999999999999
    7  = 0.9
    8  = 1.02
    9  = 1.13
    10 = 0.92
    11 = 1.13
    12 = 0.91
    13 = 1.04
    14 = 1.07
    15 = 1.08
    16 = 1.07
    17 = 1.08
}

foreach ($row in $seasonality.Keys) {
    $wsForecast.Range("L$row").Value = $seasonality[$row]
}

# --- Forecast Comparison sheet: column D (MyForecast) ---
$myForecast = @{
    5  = 6
    6  = 7
    7  = 6
    8  = 7
    9  = 7
    10 = 6
    11 = 6
    12 = 7
    13 = 6
    14 = 6
    15 = 6
    16 = 6
    17 = 5
}

foreach ($row in $myForecast.Keys) {
    $wsForecast.Range("D$row").Value = $myForecast[$row]
}

# --- Summary sheet: forecast totals (stored as text) ---
$wsSummary.Range("B9").Value  = "111"
$wsSummary.Range("B10").Value = "59"
$wsSummary.Range("B11").Value = "31"
$wsSummary.Range("B12").Value = "11"
